# Import Excel data for updating status
# Populates the sheet with an email/status table and hyperlinks each
# email address to a "mailto:" link (mirrors Excel's auto-hyperlink +
# copy/paste + per-row edit behaviour seen in the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "status"

# Row 2 - first email, hyperlinked individually
$ws.Range("A2").Value = "abc@sample.com"
$ws.Range("B2").Value = 0
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:abc@sample.com")

# Rows 3-6 - seeded by "copy/paste" of the A2 hyperlink down the column
$ws.Range("A3:A6").Value = "abc@sample.com"
$ws.Hyperlinks.Add($ws.Range("A3:A6"), "mailto:abc@sample.com", "", "", "abc@sample.com")

# Each row then gets its real email typed in, which (re)creates its own
# mailto hyperlink on top of the pasted one.
$ws.Range("A3").Value = "def@sample.com"
$ws.Range("B3").Value = 1
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:def@sample.com")

$ws.Range("A4").Value = "efg@sample.com"
$ws.Range("B4").Value = 2
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:efg@sample.com")

$ws.Range("A5").Value = "ijk@sample.com"
$ws.Range("B5").Value = 3
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:ijk@sample.com")

$ws.Range("A6").Value = "lmn@sample.com"
$ws.Range("B6").Value = 4
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:lmn@sample.com")
